$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97-122 down to 98-123.
$ws.Rows.Item(97).Insert()

# Populate the new row 97 with the new weekly price record.
$ws.Cells.Item(97, 1).Value = 4
$ws.Cells.Item(97, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(97, 3).Value = "Los Lagos"
$ws.Cells.Item(97, 4).Value = 44964
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
$ws.Cells.Item(97, 5).Value = 10
$ws.Cells.Item(97, 6).Value = 100112031
$ws.Cells.Item(97, 7).Value = "Poroto verde"
$ws.Cells.Item(97, 8).Value = "Magnum"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 50
$ws.Cells.Item(97, 11).Value = 25000
$ws.Cells.Item(97, 12).Value = 27000
$ws.Cells.Item(97, 13).Value = 26000
$ws.Cells.Item(97, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(97, 15).Value = "Región Metropolitana"
$ws.Cells.Item(97, 16).Value = 1040
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"
